$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Normalize "gender" (column C) to lowercase for data rows 2-7
$ws.Range("C2").Value = "male"
$ws.Range("C3").Value = "male"
$ws.Range("C4").Value = "male"
$ws.Range("C5").Value = "male"
$ws.Range("C6").Value = "male"
$ws.Range("C7").Value = "male"

# Normalize "location" (column D) to lowercase for data rows 2-7
$ws.Range("D2").Value = "home"
$ws.Range("D3").Value = "home"
$ws.Range("D4").Value = "home"
$ws.Range("D5").Value = "home"
$ws.Range("D6").Value = "home"
$ws.Range("D7").Value = "home"

# Normalize "activity" (column F) to lowercase for data rows 2-7 (row 2 keeps its original casing)
$ws.Range("F2").Value = "Studying"
$ws.Range("F3").Value = "studying"
$ws.Range("F4").Value = "studying"
$ws.Range("F5").Value = "studying"
$ws.Range("F6").Value = "amusing"
$ws.Range("F7").Value = "amusing"

# Normalize "mood" (column E) to lowercase for data rows 2-7
$ws.Range("E2").Value = "peaceful"
$ws.Range("E3").Value = "peaceful"
$ws.Range("E4").Value = "peaceful"
$ws.Range("E5").Value = "peaceful"
$ws.Range("E7").Value = "sad"
$ws.Range("E6").Value = "calm"

# Update the last active selection to match the final edited cell
[void]$ws.Range("F6").Select()
